$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Odd_BTTS_No (V4) 1.73 -> 1.69
$ws.Range("V4").Value = 1.69

# Row 5: Odd_Under05_FT (N5) 7.5 -> 8
$ws.Range("N5").Value = 8

# Row 5: Odd_Over25_FT (Q5) 2.25 -> 2.3
$ws.Range("Q5").Value = 2.3

# Row 5: Odd_Under25_FT (R5) 1.62 -> 1.6
$ws.Range("R5").Value = 1.6

# Row 5: Odd_BTTS_No (V5) 1.67 -> 1.63
$ws.Range("V5").Value = 1.63
